# First draft for new batch template workflow:
# - the locator strings for the new "template name"/"template description"
#   fields were captured with a stray leading space; trim them so the
#   XPath locators are clean.
# - leave the selection on the last populated row (B17) instead of the
#   stale B22 that pointed past the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "//input[@name='tname']"
$ws.Range("C11").Value = "//input[@name='tdesc']"

$ws.Range("B17").Select() | Out-Null
